$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74
$ws.Range("H74").Value = 3935.15
$ws.Range("I74").Value = 3769.9
$ws.Range("J74").Value = 4100.4
$ws.Range("K74").Value = 3769.9
$ws.Range("L74").Value = 4100.4
$ws.Range("M74").Value = -2833.9
$ws.Range("N74").Value = -5972.4
# Row 77
$ws.Range("H77").Value = 3935.15
$ws.Range("I77").Value = 3769.9
$ws.Range("J77").Value = 4100.4
$ws.Range("K77").Value = 18849.5
$ws.Range("L77").Value = 20502
$ws.Range("M77").Value = -14169.5
$ws.Range("N77").Value = -29862
# Row 86
$ws.Range("H86").Value = 2871.2
$ws.Range("I86").Value = 2561.8
$ws.Range("J86").Value = 3180.6
$ws.Range("K86").Value = 2561.8
$ws.Range("L86").Value = 3180.6
$ws.Range("M86").Value = -1438.8
$ws.Range("N86").Value = -5426.6
# Row 89
$ws.Range("H89").Value = 2871.2
$ws.Range("I89").Value = 2561.8
$ws.Range("J89").Value = 3180.6
$ws.Range("K89").Value = 12809
$ws.Range("L89").Value = 15903
$ws.Range("M89").Value = -7193
$ws.Range("N89").Value = -27135
# Row 125
$ws.Range("H125").Value = 11211576
$ws.Range("I125").Value = 581.25
$ws.Range("K125").Value = 5231.25
$ws.Range("M125").Value = -2771.25
# Row 138
$ws.Range("H138").Value = 4163468.5
$ws.Range("I138").Value = 1122817.1
$ws.Range("J138").Value = 5955281
$ws.Range("K138").Value = 3368451.3
$ws.Range("L138").Value = 17865843
$ws.Range("M138").Value = -3363311.3
$ws.Range("N138").Value = -17876123

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18771.307
$ws.Range("I32").Value = 2467.9153
$ws.Range("J32").Value = 339404.66
$ws.Range("K32").Value = 2467.9153
$ws.Range("L32").Value = 339404.66
$ws.Range("M32").Value = -2180.9153
$ws.Range("N32").Value = -339978.66
# Row 45
$ws.Range("H45").Value = 913.5833
$ws.Range("I45").Value = 866.55554
$ws.Range("K45").Value = 866.55554
$ws.Range("M45").Value = -489.55554
# Row 97
$ws.Range("H97").Value = 5599.4736
$ws.Range("I97").Value = 5599.4736
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 5599.4736
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -5103.4736
$ws.Range("N97").ClearContents() | Out-Null
# Row 102
$ws.Range("H102").Value = 1462.2222
$ws.Range("I102").Value = 1308.5714
$ws.Range("K102").Value = 1308.5714
$ws.Range("M102").Value = 313.4286
# Row 122
$ws.Range("H122").Value = 1437.0222
$ws.Range("I122").Value = 1168.4667
$ws.Range("K122").Value = 3505.4001
$ws.Range("M122").Value = -1055.4001

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1199.6666
$ws.Range("I94").Value = 1300.5
$ws.Range("J94").Value = 796.3333
$ws.Range("K94").Value = 1300.5
$ws.Range("L94").Value = 796.3333
$ws.Range("M94").Value = -849.5
$ws.Range("N94").Value = -1698.3333
# Row 134
$ws.Range("H134").Value = 2889.3809
$ws.Range("I134").Value = 1739.6875
$ws.Range("J134").Value = 6568.4
$ws.Range("K134").Value = 5219.0625
$ws.Range("L134").Value = 19705.2
$ws.Range("M134").Value = -2684.0625
$ws.Range("N134").Value = -24775.2

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 38202
$ws.Range("I2").Value = 10500
$ws.Range("J2").Value = 56670
$ws.Range("K2").Value = 10500
$ws.Range("L2").Value = 56670
$ws.Range("M2").Value = -10387
$ws.Range("N2").Value = -56896
# Row 31
$ws.Range("H31").Value = 1181.931
$ws.Range("I31").Value = 973.1852
$ws.Range("K31").Value = 973.1852
$ws.Range("M31").Value = -678.1852
# Row 34
$ws.Range("H34").Value = 1181.931
$ws.Range("I34").Value = 973.1852
$ws.Range("K34").Value = 973.1852
$ws.Range("M34").Value = -771.1852
# Row 132
$ws.Range("H132").Value = 2982.3845
$ws.Range("I132").Value = 2464.95
$ws.Range("J132").Value = 4707.1665
$ws.Range("K132").Value = 7394.849999999999
$ws.Range("L132").Value = 14121.4995
$ws.Range("M132").Value = -4864.849999999999
$ws.Range("N132").Value = -19181.4995

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1363.6818
$ws.Range("I5").Value = 735.0909
$ws.Range("K5").Value = 2205.2727
$ws.Range("M5").Value = -2093.2727
# Row 113
$ws.Range("H113").Value = 11364407
$ws.Range("I113").Value = 592.4545000000001
$ws.Range("J113").Value = 15152345
$ws.Range("K113").Value = 1777.3635
$ws.Range("L113").Value = 45457035
$ws.Range("M113").Value = 392.6364999999998
$ws.Range("N113").Value = -45461375
# Row 131
$ws.Range("H131").Value = 1330.9595
$ws.Range("I131").Value = 385.92307
$ws.Range("J131").Value = 1532.3606
$ws.Range("K131").Value = 1157.76921
$ws.Range("L131").Value = 4597.0818
$ws.Range("M131").Value = 3882.23079
$ws.Range("N131").Value = -14677.0818
# Row 135
$ws.Range("H135").Value = 1363.6818
$ws.Range("I135").Value = 735.0909
$ws.Range("K135").Value = 6615.8181
$ws.Range("M135").Value = -4080.8181

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 433.33334
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents() | Out-Null
# Row 97
$ws.Range("H97").Value = 1201
$ws.Range("I97").Value = 1201
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1201
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -705
$ws.Range("N97").ClearContents() | Out-Null
# Row 113
$ws.Range("H113").Value = 1096.25
$ws.Range("I113").Value = 1009.82355
$ws.Range("J113").Value = 1306.1428
$ws.Range("K113").Value = 1009.82355
$ws.Range("L113").Value = 1306.1428
$ws.Range("M113").Value = 1160.17645
$ws.Range("N113").Value = -5646.1428

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1549.2
$ws.Range("I16").Value = 1549.2
$ws.Range("K16").Value = 1549.2
$ws.Range("M16").Value = -1379.2
# Row 22
$ws.Range("H22").Value = 15278.714
$ws.Range("I22").Value = 433.66666
$ws.Range("J22").Value = 26412.5
$ws.Range("K22").Value = 433.66666
$ws.Range("L22").Value = 26412.5
$ws.Range("M22").Value = -138.66666
$ws.Range("N22").Value = -27002.5
# Row 27
$ws.Range("H27").Value = 15278.714
$ws.Range("I27").Value = 433.66666
$ws.Range("J27").Value = 26412.5
$ws.Range("K27").Value = 433.66666
$ws.Range("L27").Value = 26412.5
$ws.Range("M27").Value = -326.66666
$ws.Range("N27").Value = -26626.5
# Row 40
$ws.Range("H40").Value = 2196.081
$ws.Range("J40").Value = 2993.3333
$ws.Range("L40").Value = 2993.3333
$ws.Range("N40").Value = -3265.3333
# Row 55
$ws.Range("H55").Value = 399.41666
$ws.Range("I55").Value = 332.33334
$ws.Range("J55").Value = 600.6667
$ws.Range("K55").Value = 332.33334
$ws.Range("L55").Value = 600.6667
$ws.Range("M55").Value = -159.33334
$ws.Range("N55").Value = -946.6667
# Row 93
$ws.Range("H93").Value = 1180.3077
$ws.Range("I93").Value = 910
$ws.Range("K93").Value = 910
$ws.Range("M93").Value = 338
